# Fix the stale "advance to next round" picks on the Bracket sheet so that
# every winner cell actually matches one of the two animals it descended
# from (the bracket had gotten out of sync with the round-1 matchups).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bracket")

$ws.Range('D6').Value  = 'Striped Polecat'
$ws.Range('N6').Value  = 'Puffer Fish'
$ws.Range('N10').Value = 'Trapdoor Spider'
$ws.Range('E12').Value = 'Striped dolphin'
$ws.Range('D14').Value = 'Striped dolphin'
$ws.Range('G16').Value = 'Okapi'
$ws.Range('F24').Value = 'Striped Rabbit'
$ws.Range('L24').Value = 'New Caledonian Crow'
$ws.Range('N26').Value = 'New Caledonian Crow'
$ws.Range('M28').Value = 'New Caledonian Crow'
$ws.Range('D30').Value = 'Kudu'
$ws.Range('H32').Value = 'Sea Otter'
$ws.Range('I32').Value = 'Sea Otter'
$ws.Range('C35').Value = 'Shrew Mole'
$ws.Range('E44').Value = 'Mara'
$ws.Range('D46').Value = 'Mara'
$ws.Range('G48').Value = 'Sea Otter'
$ws.Range('K48').Value = 'Emperor Penguin'
$ws.Range('N50').Value = 'Bat-Eared Fox'
$ws.Range('F56').Value = 'Rock Hyrax '
$ws.Range('L56').Value = 'Greater Rhea'
$ws.Range('D58').Value = 'Bulldog Bat'
$ws.Range('M60').Value = 'Greater Rhea'
